$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "F0002"
$ws.Range("A4").Value = "F0003"
